$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Level I - Nuevo")

# Move the "Revisión Prueba / CAPM" note from Monday (A15) to Wednesday (C15),
# picking up Monday's cell formatting (centered grey style) along the way.
$ws.Range("C15").Value = "Revisión Prueba`nCAPM"
$ws.Range("A15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A15").Value = ""

# Replace the Thursday (D15) note with "Modelo Mercado"
$ws.Range("D15").Value = "Modelo Mercado`n"

# Simplify the Monday (A18) note to just "APT"
$ws.Range("A18").Value = "APT"

# Update the selected cell shown when the file was last saved
$ws.Range("A21").Select()
